$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their original text representation
# (avoid Excel auto-converting numeric-looking strings into floats/sci notation)
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D24", "D26", "D27", "D30", "D31", "D33", "D34", "D36", "D38", "D40", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = '70.585.12'
$ws.Range("E2").Value = '  +2.30%  '
$ws.Range("D3").Value = '3.814.16'
$ws.Range("E3").Value = '  +0.99%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '681.18'
$ws.Range("E5").Value = '  +8.52%  '
$ws.Range("D6").Value = '171.21'
$ws.Range("E6").Value = '  +3.63%  '
$ws.Range("D7").Value = '3.813.90'
$ws.Range("E7").Value = '  +1.03%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").Value = '0.525'
$ws.Range("E9").Value = '  +0.80%  '
$ws.Range("E10").Value = '  +1.68%  '
$ws.Range("D11").Value = '7.24'
$ws.Range("E11").Value = '  +6.61%  '
$ws.Range("D12").Value = '0.461'
$ws.Range("E12").Value = '  +0.58%  '
$ws.Range("D13").Value = '0.0000245'
$ws.Range("E13").Value = '  +0.24%  '
$ws.Range("D14").Value = '35.96'
$ws.Range("E14").Value = '  +2.00%  '
$ws.Range("D15").Value = '4.456.10'
$ws.Range("E15").Value = '  +0.99%  '
$ws.Range("D16").Value = '3.814.28'
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("D17").Value = '70.715.55'
$ws.Range("E17").Value = '  +2.42%  '
$ws.Range("D18").Value = '17.69'
$ws.Range("E18").Value = '  +0.50%  '
$ws.Range("D19").Value = '7.18'
$ws.Range("E19").Value = '  +2.05%  '
$ws.Range("E20").Value = '  +0.69%  '
$ws.Range("D21").Value = '11.27'
$ws.Range("E21").Value = '  +18.29%  '
$ws.Range("D22").Value = '477.11'
$ws.Range("E22").Value = '  +2.36%  '
$ws.Range("E23").Value = '  +1.04%  '
$ws.Range("D24").Value = '83.46'
$ws.Range("E24").Value = '  +0.58%  '
$ws.Range("E25").Value = '  -1.65%  '
$ws.Range("D26").Value = '12.27'
$ws.Range("E26").Value = '  +2.16%  '
$ws.Range("D27").Value = '10.35'
$ws.Range("E27").Value = '  +3.12%  '
$ws.Range("E28").Value = '  -1.62%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").Value = '3.964.96'
$ws.Range("E30").Value = '  +0.97%  '
$ws.Range("D31").Value = '2.91'
$ws.Range("E31").Value = '  +9.24%  '
$ws.Range("E32").Value = '  +2.39%  '
$ws.Range("D33").Value = '7.41'
$ws.Range("E33").Value = '  +3.91%  '
$ws.Range("D34").Value = '29.57'
$ws.Range("E34").Value = '  +2.93%  '
$ws.Range("E35").Value = '  +4.15%  '
$ws.Range("D36").Value = '9.13'
$ws.Range("E36").Value = '  +2.28%  '
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("D38").Value = '3.769.85'
$ws.Range("E38").Value = '  +1.13%  '
$ws.Range("E39").Value = '  +1.21%  '
$ws.Range("D40").Value = '3.39'
$ws.Range("E40").Value = '  +2.47%  '
$ws.Range("E41").Value = '  +2.29%  '
$ws.Range("E42").Value = '  -0.21%  '
$ws.Range("E44").Value = '  +11.83%  '
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("D46").Value = '46.12'
$ws.Range("E46").Value = '  +7.18%  '
$ws.Range("D47").Value = '159.94'
$ws.Range("E47").Value = '  +3.20%  '
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").Value = '1.46'
$ws.Range("E48").Value = '  +7.83%  '
$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").Value = '48.19'
$ws.Range("E49").Value = '  +3.30%  '
$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").Value = '0.000297'
$ws.Range("E50").Value = '  +9.42%  '
$ws.Range("D51").Value = '0.300'
$ws.Range("E51").Value = '  +1.51%  '
